# Generate Report for Handback
# Applies the localization-status.xlsx "handback" update:
#  - Overview sheet: status text changes from "Ready for handoff" to
#    "Handed back: in sync with en-US", and the two status columns widen.
#  - zh-cn / de-de sheets: "Latest Target File" (col I) and
#    "Latest Handback File" (col J) get populated for both data rows, the
#    "Latest Handback DateTime" (col K) is stamped, and columns C/I/J widen.

$wb = $excel.ActiveWorkbook

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7c6d00bf13bb2cb33fc88fb0a0002e0f0a7adc4/e2e/"
$file1 = "70ddb21b-ec54-4103-819b-9f7406e2b035"
$file2 = "8bd78fda-8183-49ca-a9bf-bde54cd5efb7"

# Excel's "stored" column width (the `width` attribute persisted in the
# worksheet XML) is ColumnWidth (in characters) + 5/6. Use this helper so
# we can target the persisted widths called out in the change directly.
function Set-StoredColumnWidth($col, [double]$storedWidth) {
    $col.ColumnWidth = $storedWidth - (5.0 / 6.0)
}

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

Set-StoredColumnWidth $wsOverview.Columns.Item(5) 29.9777047293527
Set-StoredColumnWidth $wsOverview.Columns.Item(6) 29.9777047293527

# ---------------------------------------------------------------------
# zh-cn sheet (row 2 = 70ddb21b..., row 3 = 8bd78fda...)
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Hyperlinks.Add($wsZh.Range("I2"), ($repoBase + $file1 + ".md"), "", "", ($file1 + ".md")) | Out-Null
$wsZh.Range("J2").Value = ($file1 + ".fbf8c160086cf3820f1c35828195f32c51160451.zh-cn.xlf")

$wsZh.Hyperlinks.Add($wsZh.Range("I3"), ($repoBase + $file2 + ".md"), "", "", ($file2 + ".md")) | Out-Null
$wsZh.Range("J3").Value = ($file2 + ".1e64ca1ae72b1ccdedeb9ad57337ae0850f7c636.zh-cn.xlf")

Set-StoredColumnWidth $wsZh.Columns.Item(3) 29.9777047293527
Set-StoredColumnWidth $wsZh.Columns.Item(9) 40
Set-StoredColumnWidth $wsZh.Columns.Item(10) 40

# ---------------------------------------------------------------------
# de-de sheet (row 2 = 70ddb21b..., row 3 = 8bd78fda...)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Hyperlinks.Add($wsDe.Range("I2"), ($repoBase + $file1 + ".md"), "", "", ($file1 + ".md")) | Out-Null
$wsDe.Range("J2").Value = ($file1 + ".fbf8c160086cf3820f1c35828195f32c51160451.de-de.xlf")
$wsDe.Range("K2").Value = "2016-08-27 19:08:34"

$wsDe.Hyperlinks.Add($wsDe.Range("I3"), ($repoBase + $file2 + ".md"), "", "", ($file2 + ".md")) | Out-Null
$wsDe.Range("J3").Value = ($file2 + ".1e64ca1ae72b1ccdedeb9ad57337ae0850f7c636.de-de.xlf")
$wsDe.Range("K3").Value = "2016-08-27 19:08:34"

Set-StoredColumnWidth $wsDe.Columns.Item(3) 29.9777047293527
Set-StoredColumnWidth $wsDe.Columns.Item(9) 40
Set-StoredColumnWidth $wsDe.Columns.Item(10) 40

# ---------------------------------------------------------------------
# zh-cn "Latest Handback DateTime" (col K) was previously the
# placeholder "0001-01-01 00:00:00" - update it now that a handback
# happened.
# ---------------------------------------------------------------------
$wsZh.Range("K2").Value = "2016-08-27 19:08:27"
$wsZh.Range("K3").Value = "2016-08-27 19:08:27"

Write-Host "Handback report generated."
